$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testcases")

# Row 6
$ws.Range("C6").Value = 'Application is installed on a desktop computer'
$ws.Range("E6").Value = '1. Copy MultiFunctionalTool_For_Desktop.zip from \\tec-share\eBX\eBX_NW\WorkArea\05._Tools\ to a preferred location<br>2. Extract the contents<br>3. Double-click on MultiFunctionalToolApplication'
$ws.Range("F6").Value = 'The application should launch successfully'

# Row 7
$ws.Range("E7").Value = '1. Copy MultiFunctionalTool_For_Laptop.zip from \\tec-share\eBX\eBX_NW\WorkArea\05._Tools\ to a preferred location<br>2. Extract the contents<br>3. Double-click on MultiFunctionalToolApplication'
$ws.Range("F7").Value = 'The application should launch successfully'

# Row 8
$ws.Range("D8").Value = 'Verify all main features are accessible from the main interface'
$ws.Range("E8").Value = '1. Launch the application<br>2. Check for Network Packet Capture feature<br>3. Check for Memory Leak Check feature<br>4. Check for Debug Log Collection feature<br>5. Check for Diagnostic Code Details feature<br>6. Check for 08 Diagnostic Code Value Access feature<br>7. Check for Protocol Configuration feature'
$ws.Range("F8").Value = 'All six main features should be accessible from the main interface'

# Row 9
$ws.Range("C9").Value = 'Application is launched'
$ws.Range("D9").Value = 'Verify Network Packet Capture - Start function'
$ws.Range("E9").Value = '1. Navigate to Network Packet Capture feature<br>2. Click on Start button'
$ws.Range("F9").Value = 'Packet capture should start successfully'

# Row 10
$ws.Range("C10").Value = 'Network Packet Capture is running'
$ws.Range("D10").Value = 'Verify Network Packet Capture - Stop function'
$ws.Range("E10").Value = '1. With packet capture running, click on Stop button'
$ws.Range("F10").Value = '1. Packet capture should stop<br>2. A .pcap file should be generated<br>3. The file should be copied to the MFP''s Shared Folder<br>4. The folder should open automatically'

# Row 11
$ws.Range("C11").Value = 'Application is launched'
$ws.Range("D11").Value = 'Verify Memory Leak Check functionality'
$ws.Range("E11").Value = '1. Navigate to Memory Leak Check feature<br>2. Select a protocol from the available options<br>3. Run the memory leak check'
$ws.Range("F11").Value = 'A protocol-specific Memory Leak Comparison Table should be displayed showing if a memory leak has occurred'

# Row 12
$ws.Range("D12").Value = 'Verify Debug Log Collection functionality'
$ws.Range("E12").Value = '1. Navigate to Debug Log Collection feature<br>2. Click on Run button'
$ws.Range("F12").Value = '1. Script execution should start<br>2. Logs should be collected<br>3. Logs should be copied to the MFP''s Shared Folder<br>4. The folder should open automatically'

# Row 13
$ws.Range("C13").Value = 'Debug Log Collection folder is empty on first attempt'
$ws.Range("D13").Value = 'Verify Debug Log Collection retry functionality'
$ws.Range("E13").Value = '1. Navigate to Debug Log Collection feature<br>2. Click on Run button<br>3. If folder is empty, run the operation again'
$ws.Range("F13").Value = 'After the second attempt, logs should be collected and visible in the MFP''s Shared Folder'

# Row 14
$ws.Range("C14").Value = 'Application is launched'
$ws.Range("D14").Value = 'Verify Diagnostic Code Details - ECC selection'
$ws.Range("E14").Value = '1. Navigate to Diagnostic Code Details feature<br>2. Select ECC option'
$ws.Range("F14").Value = 'Relevant job-specific details for ECC should be displayed'

# Row 15
$ws.Range("C15").Value = 'Application is launched'
$ws.Range("D15").Value = 'Verify Diagnostic Code Details - Network Protocols selection'
$ws.Range("E15").Value = '1. Navigate to Diagnostic Code Details feature<br>2. Select Network Protocols option'
$ws.Range("F15").Value = 'Relevant job-specific details for Network Protocols should be displayed'

# Row 16
$ws.Range("C16").Value = 'Application is launched'
$ws.Range("D16").Value = 'Verify Diagnostic Code Details - High Security Mode selection'
$ws.Range("E16").Value = '1. Navigate to Diagnostic Code Details feature<br>2. Select High Security Mode option'
$ws.Range("F16").Value = 'Relevant job-specific details for High Security Mode should be displayed'

# Row 17
$ws.Range("D17").Value = 'Verify Diagnostic Code Details - Common diagnostic codes selection'
$ws.Range("E17").Value = '1. Navigate to Diagnostic Code Details feature<br>2. Select a common diagnostic code'
$ws.Range("F17").Value = 'Relevant job-specific details for the selected diagnostic code should be displayed'

# Row 18
$ws.Range("D18").Value = 'Verify 08 Diagnostic Code Value - Get functionality'
$ws.Range("E18").Value = '1. Navigate to 08 Diagnostic Code Value Access feature<br>2. Select a diagnostic code<br>3. Click on Get button'
$ws.Range("F18").Value = 'The current value of the selected 08 diagnostic code should be displayed'

# Row 19
$ws.Range("D19").Value = 'Verify 08 Diagnostic Code Value - Set functionality'
$ws.Range("E19").Value = '1. Navigate to 08 Diagnostic Code Value Access feature<br>2. Select a diagnostic code<br>3. Enter a new value<br>4. Click on Set button'
$ws.Range("F19").Value = 'The value of the selected 08 diagnostic code should be updated successfully'

# Row 20
$ws.Range("D20").Value = 'Verify Protocol Configuration - Get functionality'
$ws.Range("E20").Value = '1. Navigate to Protocol Configuration feature<br>2. Open Protocol Selection Window<br>3. Select a protocol<br>4. Click on Get button'
$ws.Range("F20").Value = 'The current value of the selected protocol should be displayed'

# Row 21
$ws.Range("C21").Value = 'Application is launched'
$ws.Range("D21").Value = 'Verify Protocol Configuration - Set functionality'
$ws.Range("E21").Value = '1. Navigate to Protocol Configuration feature<br>2. Open Protocol Selection Window<br>3. Select a protocol<br>4. Enter a new value<br>5. Click on Set button'
$ws.Range("F21").Value = 'A message should indicate that this feature is not yet implemented'

# Row 22
$ws.Range("C22").Value = 'Application is installed and MFP is connected'
$ws.Range("D22").Value = 'Verify application performance with multiple operations'
$ws.Range("E22").Value = '1. Launch the application<br>2. Start Network Packet Capture<br>3. Open Memory Leak Check<br>4. Open Diagnostic Code Details<br>5. Stop Network Packet Capture'
$ws.Range("F22").Value = 'All operations should execute without significant delay or application crashes'
$ws.Range("H22").Value = 'Performance test'

# Row 23
$ws.Range("C23").Value = 'Application is installed'
$ws.Range("D23").Value = 'Verify application startup time'
$ws.Range("E23").Value = '1. Double-click on MultiFunctionalToolApplication<br>2. Measure the time taken for the application to fully load'
$ws.Range("F23").Value = 'Application should start within 5 seconds'
$ws.Range("H23").Value = 'Performance test'
